$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text value corrections (replace placeholder/bogus data with real data) ---
$ws.Range("B2").Value = "BBSRC-BB-V030557-1"
$ws.Range("B4").Value = "https://doi.org/10.21769/BioProtoc.126"
$ws.Range("B5").Value = "https://doi.org/10.1038/nprot.2006.232    https://doi.org/10.1038/nprot.2009.12"

# --- Fill in previously-missing measurements for the "+SUC" (LD) rows ---
$ws.Range("I12").Value = 1.1
$ws.Range("J12").Value = 1.7

$ws.Range("G13").Value = 0.2104
$ws.Range("H13").Value = 6.2
$ws.Range("I13").Value = 1.3
$ws.Range("J13").Value = 2.1

$ws.Range("G14").Value = 0.2435
$ws.Range("H14").Value = 7
$ws.Range("I14").Value = 1.2
$ws.Range("J14").Value = 1.9

$ws.Range("G15").Value = 0.3213
$ws.Range("H15").Value = 5.8
$ws.Range("I15").Value = 1.1

$ws.Range("G16").Value = 0.2135
$ws.Range("H16").Value = 4.9
$ws.Range("I16").Value = 0.8
$ws.Range("J16").Value = 2.2

$ws.Range("G17").Value = 0.292
$ws.Range("H17").Value = 5.9
$ws.Range("I17").Value = 0.9
$ws.Range("J17").Value = 2.1

# --- New "Cabinet" info row alongside the SD/LD legend (row 19) ---
$ws.Range("E19").Value = "Cabinet"
$ws.Range("F19").Value = "Percivals E-36L"
$ws.Range("E19").VerticalAlignment = -4108
$ws.Range("E19").WrapText = $true
$ws.Range("F19").VerticalAlignment = -4108
$ws.Range("F19").WrapText = $true

# --- Center-align the sample-name column (A8:A21) ---
$ws.Range("A8:A21").HorizontalAlignment = -4108

# --- New blank, centered-styled row 18 between the data rows and the legend ---
$ws.Range("A18").HorizontalAlignment = -4108

# --- Move the active selection ---
$ws.Range("D23").Select()
